# Apply updated crypto price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "71.154.96"
$ws.Cells.Item(2, 5).Value = "  +0.39%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.870.35"
$ws.Cells.Item(3, 5).Value = "  +1.72%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "699.14"
$ws.Cells.Item(5, 5).Value = "  +0.22%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "173.99"
$ws.Cells.Item(6, 5).Value = "  +0.31%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "3.868.20"
$ws.Cells.Item(7, 5).Value = "  +1.72%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.01%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +0.27%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +0.05%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "7.15"
$ws.Cells.Item(11, 5).Value = "  -5.44%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.24%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000260"
$ws.Cells.Item(13, 5).Value = "  +3.63%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "36.51"
$ws.Cells.Item(14, 5).Value = "  +0.68%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "4.521.92"
$ws.Cells.Item(15, 5).Value = "  +1.59%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "3.875.88"
$ws.Cells.Item(16, 5).Value = "  +1.96%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "71.252.28"
$ws.Cells.Item(17, 5).Value = "  +0.47%  "

# Row 18
$ws.Cells.Item(18, 2).Value = "Chainlink"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "17.73"
$ws.Cells.Item(18, 5).Value = "  -0.23%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "Polkadot"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.26"
$ws.Cells.Item(19, 5).Value = "  +0.78%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.29%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "11.18"
$ws.Cells.Item(21, 5).Value = "  -1.83%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "500.37"
$ws.Cells.Item(22, 5).Value = "  +4.35%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +1.35%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "84.96"
$ws.Cells.Item(24, 5).Value = "  +1.56%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +1.57%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "10.85"
$ws.Cells.Item(26, 5).Value = "  +5.03%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -0.70%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +0.85%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +3.08%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -0.05%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.61"

# Row 32
$ws.Cells.Item(32, 5).Value = "  -1.39%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "29.79"
$ws.Cells.Item(33, 5).Value = "  +0.68%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +1.75%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "9.27"
$ws.Cells.Item(35, 5).Value = "  +0.90%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "3.821.44"
$ws.Cells.Item(36, 5).Value = "  +1.73%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -0.51%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +2.66%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.40"
$ws.Cells.Item(39, 5).Value = "  +9.19%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -3.27%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +8.51%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +1.22%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -0.09%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "FLOKI"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.000313"
$ws.Cells.Item(45, 5).Value = "  -6.81%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "Monero"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "163.71"
$ws.Cells.Item(46, 5).Value = "  +2.08%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "48.91"
$ws.Cells.Item(47, 5).Value = "  -1.08%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +1.60%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "417.61"
$ws.Cells.Item(49, 5).Value = "  +4.48%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -2.55%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "43.78"
$ws.Cells.Item(51, 5).Value = "  -2.70%  "
